$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# --- Paragraph 2: split the single run into two runs ---
$p2 = $d.Paragraphs.Item(2)
$p2xml = '<w:p ' + $wNs + ' w14:paraId="6851725B" w14:textId="0161006A" w:rsidR="00E6682A" w:rsidRDefault="00E6682A">' `
  + '<w:r><w:t>My experiences in programming lie mainly in Java, C and Python</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve"> where I have created simple side-scrolling games, interactive environments and sorting mechanisms that help users organize their file space. </w:t></w:r>' `
  + '</w:p>'
$p2.Range.InsertXML($p2xml)

# --- Paragraph 3: split the single run into six runs, keep the bookmark ---
$p3 = $d.Paragraphs.Item(3)
$p3xml = '<w:p ' + $wNs + ' w14:paraId="01599D1B" w14:textId="04D9AC2F" w:rsidR="00E6682A" w:rsidRDefault="00E6682A">' `
  + '<w:r><w:t>I am a game designer whose interests lie in creating games with strong narratives and transformative experiences. My goals in creating these games are to evoke a deep</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve"> and meaningful</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve"> feeling </w:t></w:r>' `
  + '<w:r><w:t>from</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve"> the player</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve"> to create experiences that they can share with others. </w:t></w:r>' `
  + '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' `
  + '</w:p>'
$p3.Range.InsertXML($p3xml)
